# Commit: "css: add CSS values and units and sizing items"
#
# Adds 7 new word/phonetic-symbol pairs to the "word" sheet and 2 new
# phrase entries (with Chinese translation + example sentence) to the
# "phrase" sheet, then leaves the workbook with the "phrase" sheet active.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("word")
$ws2 = $wb.Worksheets.Item("phrase")

# ---------------------------------------------------------------------
# 1. "word" sheet - append new word / phonetic-symbol rows (59-65)
# ---------------------------------------------------------------------
$ws1.Range("A59").Value = "opaque"
$ws1.Range("B59").Value = "/oʊˈpeɪk/"

$ws1.Range("A60").Value = "opacity"
$ws1.Range("B60").Value = "/oʊˈpæsəti/"

$ws1.Range("A61").Value = "finite"
$ws1.Range("B61").Value = "/ˈfaɪnaɪt/"

$ws1.Range("A62").Value = "infinite"
$ws1.Range("B62").Value = "/ˈɪnfɪnət/"

$ws1.Range("A63").Value = "align"
$ws1.Range("B63").Value = "/əˈlaɪn/"

$ws1.Range("A64").Value = "intrinsic"
$ws1.Range("B64").Value = "/ɪnˈtrɪnzɪk/"

$ws1.Range("A65").Value = "extrinsic"
$ws1.Range("B65").Value = "/ɪksˈtrɪnzɪk/"

# restore the "word" sheet's own last selection
$ws1.Range("C71").Select()

# ---------------------------------------------------------------------
# 2. "phrase" sheet - append new phrase rows (73-74)
# ---------------------------------------------------------------------
$ws2.Range("A73").Value = "over and over"
$ws2.Range("B73").Value = "反复、再三、一次又一次"
$ws2.Range("D73").Value = "In programming, a function is a piece of code that does a specific task. Functions are useful because you can write code once then reuse it many times instead of writing the same logic over and over."

$ws2.Range("A74").Value = "scale down"
$ws2.Range("B74").Value = "按比例缩小"
$ws2.Range("D74").Value = "A common use of max-width is to cause images to scale down if there is not enough space to display them at their intrinsic width, while making sure they don't become larger than that width."

# make "phrase" the active sheet/tab with the given selection in the
# frozen (bottomLeft) pane, matching the final workbook view state
$ws2.Activate()
$ws2.Range("B79").Select()

Write-Output "Added 7 word entries and 2 phrase entries."
